$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 935.0903682874737
$ws.Range("D2").Value = 4.6754518414373685
$ws.Range("E2").Value = 6.493683113107456

$ws.Range("B3").Value = 701.3177762156054
$ws.Range("D3").Value = 3.506588881078027
$ws.Range("E3").Value = 4.870262334830593

$ws.Range("B4").Value = 124.67871577166318
$ws.Range("D4").Value = 0.6233935788583159
$ws.Range("E4").Value = 0.8658244150809943

$ws.Range("B5").Value = 1761.086860274742
$ws.Range("D5").Value = 8.80543430137371
$ws.Range("E5").Value = 12.229769863019042

$ws.Range("B7").Value = 258.08494164734276
$ws.Range("D7").Value = 1.2904247082367137
$ws.Range("E7").Value = 1.7922565392176582

$ws.Range("B8").Value = 43.014156941223796
$ws.Range("D8").Value = 0.21507078470611898
$ws.Range("E8").Value = 0.29870942320294297

$ws.Range("B9").Value = 301.09909858856656
$ws.Range("D9").Value = 1.5054954929428328
$ws.Range("E9").Value = 2.090965962420601

$ws.Range("B11").Value = 276.3110414475657
$ws.Range("C11").Value = 128.4744656626086
$ws.Range("D11").Value = 1.3815552072378283
$ws.Range("E11").Value = 1.9188266767192061

$ws.Range("B13").Value = 171.25623456985784
$ws.Range("C13").Value = 79.62784661983216
$ws.Range("D13").Value = 0.8562811728492892
$ws.Range("E13").Value = 1.1892794067351238

$ws.Range("B18").Value = 108.74770895185972
$ws.Range("C18").Value = 50.56368260359341
$ws.Range("D18").Value = 0.5437385447592986
$ws.Range("E18").Value = 0.7551924232768036

$ws.Range("B20").Value = 395.5332959134297
$ws.Range("C20").Value = 183.9084264531334
$ws.Range("D20").Value = 1.9776664795671484
$ws.Range("E20").Value = 2.7467589993988173

$ws.Range("B21").Value = 282.6580047461431
$ws.Range("C21").Value = 131.42557001983226
$ws.Range("D21").Value = 1.4132900237307155
$ws.Range("E21").Value = 1.9629028107371047

$ws.Range("B22").Value = 760.2117748158623
$ws.Range("C22").Value = 353.4704984940865
$ws.Range("D22").Value = 3.8010588740793114
$ws.Range("E22").Value = 5.279248436221266

$ws.Range("B25").Value = 3207.456484078597
$ws.Range("C25").Value = 1491.3492264704337
$ws.Range("D25").Value = 16.037282420392984
$ws.Range("E25").Value = 22.274003361656924

$ws.Range("B27").Value = 1446.3696238038542
$ws.Range("C27").Value = 672.5086467602885
$ws.Range("D27").Value = 7.231848119019271
$ws.Range("E27").Value = 10.044233498637876
